$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new task row (row 28) following the existing "Tareas" table layout:
# A=Tarea, B=Horas estimadas, C=Horas reales, D=Autor, E=Fecha
$ws.Range("A28").Value = "Corregir test"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "Clara"
$ws.Range("E28").Value = Get-Date -Year 2016 -Month 10 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("E28").NumberFormat = "mm/dd/yyyy"

# Match the formatting used on the row above (style carried from A27/E27)
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection/view to mirror the edit location
$ws.Range("B28").Select()
